$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "1.00", "0.600").
# Force text format before assigning so Excel does not silently coerce
# these into numbers and drop formatting like trailing zeros or thousands dots.
$dCells = @("D2","D3","D5","D6","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D43","D44","D45","D47","D48","D50","D51")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range('D2').Value = '83.812.42'
$ws.Range('E2').Value = '  +5.15%  '
$ws.Range('D3').Value = '3.225.13'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '218.43'
$ws.Range('E5').Value = '  +2.82%  '
$ws.Range('D6').Value = '623.67'
$ws.Range('E6').Value = '  -2.12%  '
$ws.Range('E7').Value = '  +25.20%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '0.589'
$ws.Range('E9').Value = '  -2.29%  '
$ws.Range('D10').Value = '3.221.65'
$ws.Range('E10').Value = '  +0.53%  '
$ws.Range('D11').Value = '0.600'
$ws.Range('E11').Value = '  +0.99%  '
$ws.Range('D12').Value = '0.0000275'
$ws.Range('E12').Value = '  +5.74%  '
$ws.Range('D13').Value = '0.165'
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.833.50'
$ws.Range('E14').Value = '  +1.03%  '
$ws.Range('B15').Value = 'Toncoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D15').Value = '5.37'
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('D16').Value = '32.54'
$ws.Range('E16').Value = '  +0.79%  '
$ws.Range('D17').Value = '83.252.37'
$ws.Range('E17').Value = '  +4.67%  '
$ws.Range('D18').Value = '3.234.48'
$ws.Range('E18').Value = '  +0.97%  '
$ws.Range('D19').Value = '3.21'
$ws.Range('E19').Value = '  +7.08%  '
$ws.Range('D20').Value = '14.31'
$ws.Range('E20').Value = '  -1.45%  '
$ws.Range('D21').Value = '445.92'
$ws.Range('E21').Value = '  +1.37%  '
$ws.Range('D22').Value = '9.03'
$ws.Range('E22').Value = '  -3.46%  '
$ws.Range('D23').Value = '5.20'
$ws.Range('E23').Value = '  -0.25%  '
$ws.Range('D24').Value = '7.43'
$ws.Range('E24').Value = '  +6.44%  '
$ws.Range('D25').Value = '5.19'
$ws.Range('E25').Value = '  +7.22%  '
$ws.Range('D26').Value = '11.87'
$ws.Range('E26').Value = '  +9.14%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '3.394.41'
$ws.Range('E27').Value = '  +0.77%  '
$ws.Range('B28').Value = 'Litecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D28').Value = '78.30'
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').Value = '0.0000123'
$ws.Range('E30').Value = '  -1.53%  '
$ws.Range('D31').Value = '9.12'
$ws.Range('E31').Value = '  -0.83%  '
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.27%  '
$ws.Range('D33').Value = '568.63'
$ws.Range('E33').Value = '  +1.10%  '
$ws.Range('D34').Value = '0.155'
$ws.Range('E34').Value = '  +27.27%  '
$ws.Range('D35').Value = '1.48'
$ws.Range('E35').Value = '  -3.27%  '
$ws.Range('D36').Value = '0.153'
$ws.Range('E36').Value = '  -2.95%  '
$ws.Range('D37').Value = '1.99'
$ws.Range('E37').Value = '  -2.60%  '
$ws.Range('D38').Value = '23.08'
$ws.Range('E38').Value = '  -0.29%  '
$ws.Range('D39').Value = '6.19'
$ws.Range('E39').Value = '  +9.05%  '
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('D41').Value = '0.407'
$ws.Range('E41').Value = '  -1.85%  '
$ws.Range('E42').Value = '  +3.31%  '
$ws.Range('D43').Value = '2.03'
$ws.Range('E43').Value = '  +10.54%  '
$ws.Range('D44').Value = '3.03'
$ws.Range('E44').Value = '  +12.10%  '
$ws.Range('D45').Value = '160.45'
$ws.Range('E45').Value = '  -1.80%  '
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('D47').Value = '188.05'
$ws.Range('E47').Value = '  -2.24%  '
$ws.Range('D48').Value = '44.90'
$ws.Range('E48').Value = '  +3.47%  '
$ws.Range('E49').Value = '  -2.12%  '
$ws.Range('D50').Value = '0.775'
$ws.Range('E50').Value = '  -3.31%  '
$ws.Range('D51').Value = '25.83'
$ws.Range('E51').Value = '  +0.52%  '

foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "General" }
